# ---------------------------------------------------------------------------
# Applies the tracked changes to Template_Presentation_Students.pptx:
#   1. Re-caches the "datetimeFigureOut" date placeholder on the slide
#      master and all five slide layouts from 6/6/2024 -> 6/13/2024.
#   2. On the last slide (slide 10), removes the two small decorative
#      corner-square shapes ("object 3" and "object 5").
#   3. Adds a new "TextBox 3" shape under the body copy on slide 10 with
#      the project's GitHub link.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Update the cached date field text on master + every layout --------

$newDate = "6/13/2024"

function Update-DatePlaceholder {
    param($shapes, [string]$text)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $text
        }
    }
}

$master = $p.Slides.Item(1).Master
Update-DatePlaceholder $master.Shapes $newDate

# Master.CustomLayouts.Item(n) is unreliable for n > 1 in this host, so we
# reach every layout by briefly instantiating a slide against it (via the
# classic ppLayout-index overload of Slides.Add) and removing that scratch
# slide again once its layout has been edited.
$layoutCount = 5
for ($layoutNum = 1; $layoutNum -le $layoutCount; $layoutNum++) {
    $scratch = $p.Slides.Add($p.Slides.Count + 1, $layoutNum)
    Update-DatePlaceholder $scratch.CustomLayout.Shapes $newDate
    $scratch.Delete()
}

# --- 2. Remove the two decorative square shapes from slide 10 -------------

$lastSlide = $p.Slides.Item($p.Slides.Count)
$lastSlide.Shapes.Item("object 3").Delete()
$lastSlide.Shapes.Item("object 5").Delete()

# --- 3. Add the new "Project Link" textbox to slide 10 --------------------

$emuPerPoint = 12700.0
$leftPt = 382718 / $emuPerPoint
$topPt = 4724400 / $emuPerPoint
$widthPt = 11468524 / $emuPerPoint
$heightPt = 523220 / $emuPerPoint

$linkBox = $lastSlide.Shapes.AddTextbox(1, $leftPt, $topPt, $widthPt, $heightPt)
$linkBox.Fill.Visible = $false
$linkBox.TextFrame.WordWrap = $false

$tr = $linkBox.TextFrame.TextRange
$tr.Text = "Project Link: https://github.com/Silent-killer-from-vizag/Cyber_Project.git"
$tr.Font.Name = "Bell MT"
$tr.Font.Size = 28

$linkBox.TextFrame.AutoSize = 1
